# Remove the "总市值2 / [交易日期] 最新收盘日 / [单位] 元" market-cap column
# (column C) and the trailing "数据来源：Wind资讯" footnote row (row 31) while
# the team tries different parameters for the excess-return calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop any column-level formatting remnants on C first so deleting it does not
# leave a stray <col> entry behind, then remove the column entirely (shifts
# nothing left of it, just collapses C out of existence).
$ws.Columns("C").ClearFormats()
$ws.Columns("C").Delete()

# Remove the footer/source note row that lived below the data (was row 31).
$ws.Rows(31).Delete()

# Row 1 only needed its extra height to fit the wrapped market-cap header that
# lived in column C; now that it's gone, let the row shrink back down.
$ws.Rows(1).AutoFit()

# Give the two remaining columns their own explicit widths (closest reachable
# to the authored 19.44140625 / 14.21875 character widths).
$ws.Columns("A").ColumnWidth = 18.61
$ws.Columns("B").ColumnWidth = 13.39

# Restore the cursor to where the author left it after the edit.
$ws.Range("D11").Select()
